$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows whose District (column G) text is exactly "Kalaburagi" and must be
# updated to the official name "Kalaburagi (Gulbarga)".
$rowsToUpdateDistrict = @(3,6,7,8,10,11,12,13,14,15,16,17,19,21,23,24,25,26,27,30,31,32,33,35,37,41,43,44,45,46,49,51,52,53,54,55,56,57,58)

foreach ($r in $rowsToUpdateDistrict) {
    $cell = $ws.Cells.Item($r, 7)   # column G
    if ($cell.Text -eq "Kalaburagi") {
        $cell.Value = "Kalaburagi (Gulbarga)"
    }
}

# Rows whose Address (column F) is an empty inline string cell that should
# be removed entirely (not merely blanked).
$rowsToClearAddress = @(4,5,9,13,18,21,22,34,39,40,42,48)

foreach ($r in $rowsToClearAddress) {
    $cell = $ws.Cells.Item($r, 6)   # column F
    $cell.ClearContents()
}
